# Apply updated cryptocurrency price/volume data (refresh scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.041.04"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "1.642.38"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.48%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.65"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0636"
$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("E9").Value = "  -1.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.44"
$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  +0.09%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.25"
$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.659.90"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.541"
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.23"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").Value = "0.0₃0758"
$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").Value = "26.049.23"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "193.74"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.33"
$ws.Range("E20").Value = "  -0.82%  "

$ws.Range("E21").Value = "  -2.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("E22").Value = "  -2.03%  "

$ws.Range("E23").Value = "  +1.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.13"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  +0.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.84"
$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.47"
$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("E30").Value = "  -3.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.26"
$ws.Range("E31").Value = "  +0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.53"
$ws.Range("E33").Value = "  -0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.46"
$ws.Range("E34").Value = "  +0.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.899"
$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("D36").Value = "1.128.85"
$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.47"
$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.533"
$ws.Range("E38").Value = "  -2.49%  "

$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.41"
$ws.Range("E40").Value = "  -0.96%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.796"
$ws.Range("E41").Value = "  -0.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.32"
$ws.Range("E42").Value = "  -3.08%  "

$ws.Range("D43").Value = "0.0₆0115"
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "56.18"
$ws.Range("E44").Value = "  -0.83%  "

$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.74"
$ws.Range("E47").Value = "  +0.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.415"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.01"
$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("E50").Value = "  -2.81%  "

$ws.Range("E51").Value = "  +0.93%  "
